# Auto-generated Excel COM-interop script applying the Ravana_Profits market-price refresh diff.
# For each affected row/cell: update numeric values, delete cells that are removed entirely,
# and add cells that are newly introduced (e.g. M97, N86, N89).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 1082.4286
$ws.Range("I2").Value = 1012.8333
$ws.Range("K2").Value = 1012.8333
$ws.Range("M2").Value = -899.8333

# Row 6
$ws.Range("H6").Value = 378
$ws.Range("I6").Value = 378
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 1134
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -1022
$ws.Range("N6").ClearContents()

# Row 12
$ws.Range("H12").Value = 679.8333
$ws.Range("J12").Value = 793
$ws.Range("L12").Value = 793
$ws.Range("N12").Value = -1133

# Row 18
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()

# Row 21
$ws.Range("H21").Value = 30000
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

# Row 23
$ws.Range("H23").Value = 30000
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

# Row 40
$ws.Range("H40").Value = 1600.3334
$ws.Range("J40").Value = 1400
$ws.Range("L40").Value = 1400
$ws.Range("N40").Value = -1750

# Row 135
$ws.Range("H135").Value = 869.04346
$ws.Range("I135").Value = 761.381
$ws.Range("K135").Value = 6852.429
$ws.Range("M135").Value = -4317.429

# Row 137
$ws.Range("H137").Value = 3693.05
$ws.Range("I137").Value = 2491.1428
$ws.Range("K137").Value = 7473.428400000001
$ws.Range("M137").Value = -4923.428400000001

# Row 138
$ws.Range("H138").Value = 4148.1055
$ws.Range("J138").Value = 5942.619
$ws.Range("L138").Value = 17827.857
$ws.Range("N138").Value = -28107.857

# Row 141
$ws.Range("H141").Value = 2990.6316
$ws.Range("I141").Value = 2628.4
$ws.Range("K141").Value = 7885.200000000001
$ws.Range("M141").Value = -2705.200000000001

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5095.3057
$ws.Range("I32").Value = 4946.394
$ws.Range("J32").Value = 6733.3335
$ws.Range("K32").Value = 4946.394
$ws.Range("L32").Value = 6733.3335
$ws.Range("M32").Value = -4659.394
$ws.Range("N32").Value = -7307.3335

# Row 61
$ws.Range("H61").Value = 1087.5294
$ws.Range("I61").Value = 1133.4667
$ws.Range("J61").Value = 743
$ws.Range("K61").Value = 1133.4667
$ws.Range("L61").Value = 743
$ws.Range("M61").Value = -921.4666999999999
$ws.Range("N61").Value = -1167

# Row 74
$ws.Range("H74").Value = 2152.0588
$ws.Range("I74").Value = 2189.75
$ws.Range("K74").Value = 2189.75
$ws.Range("M74").Value = -1315.75

# Row 77
$ws.Range("H77").Value = 2152.0588
$ws.Range("I77").Value = 2189.75
$ws.Range("K77").Value = 10948.75
$ws.Range("M77").Value = -6580.75

# Row 102
$ws.Range("H102").Value = 2103.4285
$ws.Range("I102").Value = 2103.4285
$ws.Range("K102").Value = 2103.4285
$ws.Range("M102").Value = -481.4285

# Row 122
$ws.Range("H122").Value = 2247.6
$ws.Range("I122").Value = 2247.6
$ws.Range("K122").Value = 6742.799999999999
$ws.Range("M122").Value = -4292.799999999999

# Row 136
$ws.Range("H136").Value = 1087.5294
$ws.Range("I136").Value = 1133.4667
$ws.Range("J136").Value = 743
$ws.Range("K136").Value = 3400.4001
$ws.Range("L136").Value = 2229
$ws.Range("M136").Value = -850.4000999999998
$ws.Range("N136").Value = -7329

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 4000
$ws.Range("J86").Value = 4000
$ws.Range("L86").Value = 4000
$ws.Range("N86").Value = -6246

# Row 89
$ws.Range("H89").Value = 4000
$ws.Range("J89").Value = 4000
$ws.Range("L89").Value = 20000
$ws.Range("N89").Value = -31232

# Row 105
$ws.Range("H105").Value = 2687.7334
$ws.Range("I105").Value = 2529.7144
$ws.Range("K105").Value = 2529.7144
$ws.Range("M105").Value = -782.7143999999998

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2554.1428
$ws.Range("I31").Value = 2616.2
$ws.Range("K31").Value = 2616.2
$ws.Range("M31").Value = -2321.2

# Row 34
$ws.Range("H34").Value = 2554.1428
$ws.Range("I34").Value = 2616.2
$ws.Range("K34").Value = 2616.2
$ws.Range("M34").Value = -2414.2

# Row 132
$ws.Range("H132").Value = 3094.4167
$ws.Range("I132").Value = 2921.2727
$ws.Range("K132").Value = 8763.8181
$ws.Range("M132").Value = -6233.8181

# Row 134
$ws.Range("H134").Value = 3117.2
$ws.Range("I134").Value = 3117.2
$ws.Range("K134").Value = 9351.599999999999
$ws.Range("M134").Value = -6816.599999999999

$ws = $wb.Worksheets.Item("CUL")
# Row 109
$ws.Range("H109").Value = 3921.3076
$ws.Range("I109").Value = 992.3333
$ws.Range("K109").Value = 2976.9999
$ws.Range("M109").Value = -1936.9999

# Row 122
$ws.Range("H122").Value = 112311.664
$ws.Range("I122").Value = 903.5
$ws.Range("J122").Value = 144142.58
$ws.Range("K122").Value = 8131.5
$ws.Range("L122").Value = 1297283.22
$ws.Range("M122").Value = -5681.5
$ws.Range("N122").Value = -1302183.22

$ws = $wb.Worksheets.Item("GSM")
# Row 22
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()

# Row 70
$ws.Range("H70").Value = 13237
$ws.Range("I70").Value = 12832.667
$ws.Range("K70").Value = 12832.667
$ws.Range("M70").Value = -12562.667

# Row 73
$ws.Range("H73").Value = 13237
$ws.Range("I73").Value = 12832.667
$ws.Range("K73").Value = 12832.667
$ws.Range("M73").Value = -11896.667

# Row 132
$ws.Range("H132").Value = 1325.174
$ws.Range("I132").Value = 746.73334
$ws.Range("J132").Value = 2409.75
$ws.Range("K132").Value = 2240.20002
$ws.Range("L132").Value = 7229.25
$ws.Range("M132").Value = 289.7999799999998
$ws.Range("N132").Value = -12289.25

$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 2403.375
$ws.Range("I93").Value = 2403.375
$ws.Range("K93").Value = 2403.375
$ws.Range("M93").Value = -1155.375

# Row 122
$ws.Range("H122").Value = 2951.1667
$ws.Range("I122").Value = 2951.1667
$ws.Range("K122").Value = 8853.500100000001
$ws.Range("M122").Value = -6403.500100000001

# Row 136
$ws.Range("H136").Value = 3017.8147
$ws.Range("I136").Value = 2907.1667
$ws.Range("K136").Value = 8721.500100000001
$ws.Range("M136").Value = -6171.500100000001

$ws = $wb.Worksheets.Item("WVR")
# Row 97
$ws.Range("H97").Value = 20260
$ws.Range("I97").Value = 520
$ws.Range("J97").Value = 40000
$ws.Range("K97").Value = 520
$ws.Range("L97").Value = 40000
$ws.Range("M97").Value = 471
$ws.Range("N97").Value = -41982

# Row 136
$ws.Range("H136").Value = 1278.4286
$ws.Range("I136").Value = 951.88
$ws.Range("K136").Value = 2855.64
$ws.Range("M136").Value = -305.6399999999999
